# modelo_calendario.xlsx update:
#  - disambiguate the three sample "Reunião Teste" rows so each event has
#    its own subject (Reunião Teste 1/2/3)
#  - add a blank, underlined row below the sample data as a place to start
#    typing new rows, and select it
#  - set the page to A4/portrait for printing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give each sample meeting row a unique subject.
$ws.Range("A2").Value = "Reunião Teste 1"
$ws.Range("A3").Value = "Reunião Teste 2"
$ws.Range("A4").Value = "Reunião Teste 3"

# New empty row right after the sample data, underlined and selected -
# ready for the next entry.
$xlUnderlineStyleSingle = 2
$ws.Range("A5").Font.Underline = $xlUnderlineStyleSingle
[void]$ws.Range("A5").Select()

# Page setup for printing.
$xlPortrait = 1
$ws.PageSetup.PaperSize = 9        # xlPaperA4
$ws.PageSetup.Orientation = $xlPortrait
